$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recalculated salary index values for rows 2-18 (columns B:H)
$ws.Cells.Item(2, 2).Value = 11.60260019864755
$ws.Cells.Item(2, 3).Value = 12.33368594715992
$ws.Cells.Item(2, 4).Value = 11.967846525
$ws.Cells.Item(2, 5).Value = 11.18011887873801
$ws.Cells.Item(2, 6).Value = 11.25746539402369
$ws.Cells.Item(2, 7).Value = 11.31191657662111
$ws.Cells.Item(2, 8).Value = 11.0948433886587

$ws.Cells.Item(3, 2).Value = 12.59231637119373
$ws.Cells.Item(3, 3).Value = 12.94189373991007
$ws.Cells.Item(3, 4).Value = 12.6781539651283
$ws.Cells.Item(3, 5).Value = 11.83621986310745
$ws.Cells.Item(3, 6).Value = 12.03783229194669
$ws.Cells.Item(3, 7).Value = 12.20090987489004
$ws.Cells.Item(3, 8).Value = 11.79205809379728

$ws.Cells.Item(4, 2).Value = 13.4667622838789
$ws.Cells.Item(4, 3).Value = 14.13864366746315
$ws.Cells.Item(4, 4).Value = 13.799922375
$ws.Cells.Item(4, 5).Value = 12.7721133334717
$ws.Cells.Item(4, 6).Value = 12.93852660929048
$ws.Cells.Item(4, 7).Value = 13.13300537461461
$ws.Cells.Item(4, 8).Value = 12.36535599334476

$ws.Cells.Item(5, 2).Value = 11.77925334390886
$ws.Cells.Item(5, 3).Value = 12.7458725529953
$ws.Cells.Item(5, 4).Value = 12.2383543236286
$ws.Cells.Item(5, 5).Value = 11.46529868898703
$ws.Cells.Item(5, 6).Value = 11.75750587755254
$ws.Cells.Item(5, 7).Value = 12.00597498558663
$ws.Cells.Item(5, 8).Value = 11.84818824381014

$ws.Cells.Item(6, 2).Value = 11.04136413156123
$ws.Cells.Item(6, 3).Value = 11.68963463147117
$ws.Cells.Item(6, 4).Value = 11.3585287154406
$ws.Cells.Item(6, 5).Value = 10.77479266543566
$ws.Cells.Item(6, 6).Value = 10.64941057978176
$ws.Cells.Item(6, 7).Value = 10.74454439942906
$ws.Cells.Item(6, 8).Value = 10.54267422872699

$ws.Cells.Item(7, 2).Value = 12.52136834631571
$ws.Cells.Item(7, 3).Value = 13.03817355621505
$ws.Cells.Item(7, 4).Value = 12.42311535259613
$ws.Cells.Item(7, 5).Value = 11.87366646679993
$ws.Cells.Item(7, 6).Value = 11.87617869085521
$ws.Cells.Item(7, 7).Value = 11.98882425984096
$ws.Cells.Item(7, 8).Value = 11.42292719915424

$ws.Cells.Item(8, 2).Value = 11.92937655860349
$ws.Cells.Item(8, 3).Value = 12.46490140278902
$ws.Cells.Item(8, 4).Value = 12.156855975
$ws.Cells.Item(8, 5).Value = 11.40461087454539
$ws.Cells.Item(8, 6).Value = 11.5807511907452
$ws.Cells.Item(8, 7).Value = 11.71345400159246
$ws.Cells.Item(8, 8).Value = 11.44121022945437

$ws.Cells.Item(9, 2).Value = 11.58063773081721
$ws.Cells.Item(9, 3).Value = 12.16988072206184
$ws.Cells.Item(9, 4).Value = 11.93286706905589
$ws.Cells.Item(9, 5).Value = 11.03170125288324
$ws.Cells.Item(9, 6).Value = 11.07281577020706
$ws.Cells.Item(9, 7).Value = 11.21224091020639
$ws.Cells.Item(9, 8).Value = 10.96747645961166

$ws.Cells.Item(10, 2).Value = 14.01261739635469
$ws.Cells.Item(10, 3).Value = 14.74228291531547
$ws.Cells.Item(10, 4).Value = 14.41524198729834
$ws.Cells.Item(10, 5).Value = 13.60834283068732
$ws.Cells.Item(10, 6).Value = 13.86919613336777
$ws.Cells.Item(10, 7).Value = 14.0061547047694
$ws.Cells.Item(10, 8).Value = 13.65358531443892

$ws.Cells.Item(11, 2).Value = 11.96245767122055
$ws.Cells.Item(11, 3).Value = 12.52913047146017
$ws.Cells.Item(11, 4).Value = 12.29499865416779
$ws.Cells.Item(11, 5).Value = 11.52187616689947
$ws.Cells.Item(11, 6).Value = 11.77064616174712
$ws.Cells.Item(11, 7).Value = 11.96075867871859
$ws.Cells.Item(11, 8).Value = 11.57075624929744

$ws.Cells.Item(12, 2).Value = 10.89174287932444
$ws.Cells.Item(12, 3).Value = 11.18313488420164
$ws.Cells.Item(12, 4).Value = 11.14584952641588
$ws.Cells.Item(12, 5).Value = 10.41962370921107
$ws.Cells.Item(12, 6).Value = 10.48201131607158
$ws.Cells.Item(12, 7).Value = 10.8828572561911
$ws.Cells.Item(12, 8).Value = 10.67089705020983

$ws.Cells.Item(13, 2).Value = 11.95100758960922
$ws.Cells.Item(13, 3).Value = 12.52890557923256
$ws.Cells.Item(13, 4).Value = 12.27911615
$ws.Cells.Item(13, 5).Value = 11.40159626372215
$ws.Cells.Item(13, 6).Value = 11.49929374661104
$ws.Cells.Item(13, 7).Value = 11.81589671469163
$ws.Cells.Item(13, 8).Value = 11.57171436027681

$ws.Cells.Item(14, 2).Value = 14.89118209854162
$ws.Cells.Item(14, 3).Value = 15.95645235579862
$ws.Cells.Item(14, 4).Value = 15.45918711734407
$ws.Cells.Item(14, 5).Value = 14.64623406208957
$ws.Cells.Item(14, 6).Value = 15.08291186356561
$ws.Cells.Item(14, 7).Value = 15.34078644749878
$ws.Cells.Item(14, 8).Value = 14.79103938427232

$ws.Cells.Item(15, 2).Value = 11.30672076401736
$ws.Cells.Item(15, 3).Value = 11.95675404627403
$ws.Cells.Item(15, 4).Value = 11.76079419933817
$ws.Cells.Item(15, 5).Value = 11.00712111789826
$ws.Cells.Item(15, 6).Value = 11.25576283749813
$ws.Cells.Item(15, 7).Value = 11.58549535169647
$ws.Cells.Item(15, 8).Value = 11.25198932239388

$ws.Cells.Item(16, 2).Value = 14.62933733082327
$ws.Cells.Item(16, 3).Value = 15.42078682931872
$ws.Cells.Item(16, 4).Value = 14.69858752382294
$ws.Cells.Item(16, 5).Value = 13.96784545364735
$ws.Cells.Item(16, 6).Value = 14.12999791611706
$ws.Cells.Item(16, 7).Value = 14.3395125069415
$ws.Cells.Item(16, 8).Value = 13.79513149636294

$ws.Cells.Item(17, 2).Value = 16.29310739315033
$ws.Cells.Item(17, 3).Value = 17.19013484234888
$ws.Cells.Item(17, 4).Value = 16.29247682706402
$ws.Cells.Item(17, 5).Value = 15.4326069921872
$ws.Cells.Item(17, 6).Value = 15.58695838087561
$ws.Cells.Item(17, 7).Value = 15.56313738617234
$ws.Cells.Item(17, 8).Value = 14.98580102321255

$ws.Cells.Item(18, 2).Value = 12.09773894536583
$ws.Cells.Item(18, 3).Value = 12.61083418147217
$ws.Cells.Item(18, 4).Value = 12.29918200136333
$ws.Cells.Item(18, 5).Value = 11.43500496148252
$ws.Cells.Item(18, 6).Value = 11.72513991261468
$ws.Cells.Item(18, 7).Value = 11.88401368990424
$ws.Cells.Item(18, 8).Value = 11.67685402695732

# Row 19 ("Total Nacional") has been removed entirely
$ws.Rows.Item(19).Delete()

# Update the sheet dimension to reflect the new used range
$ws.Range("A1:H18").Select() | Out-Null
